$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.166.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.418.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.34%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'555.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.41%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'159.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.506"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.52%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.162"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +6.74%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.73%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -5.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'68.038.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.40%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.856.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.71%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +2.19%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'22.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.92%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.413.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.98%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'10.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.62%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'331.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.14%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.39%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'SuiNetwork"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'1.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.99%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'Dai"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.04%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'66.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.62%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.537.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.55%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.13%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0808"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.89%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.78%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.11%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'422.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.30%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.95%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.18%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.02%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'17.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.45%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.38%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.51%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.21%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'132.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -1.32%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.91%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.62%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.555"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0913"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.74%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.20%  "
$ws.Range("E51").Style = "Normal"
